$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new rows right after the current last data row (24) ---
$ws.Rows("25:26").Insert()

# --- 2. Populate new row 25 with PATRICIA's data (copy format+values from row 23, the "normal" style row) ---
$ws.Range("B23:J23").Copy($ws.Range("B25:J25"))

# --- 3. Populate new row 26 with CLEIVER's data (copy format+values from row 24, the "last row / special border" style row) ---
$ws.Range("B24:J24").Copy($ws.Range("B26:J26"))

# --- 4. Row 24 is no longer the last row, so it should lose the special bottom-border style ---
#        and pick up the "normal" formatting used by the rest of the table (like row 23).
$ws.Range("B23:J23").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)

# --- 5. New period "2509" applies to rows 25 and 26 (was 2508, copied from rows 23/24) ---
$ws.Range("E25").Value = "2509"
$ws.Range("E26").Value = "2509"

# --- 6. Update the summary figures: Valor Mora total and Cant. Periodos ---
$ws.Range("E11").Value = 490724
$ws.Range("F13").Value = 8
